# Update "想去人数" (F column) values across sheets per upstream data refresh (commit 456a3b4)
$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("展览")
$ws.Range("F5").Value = 5769
$ws.Range("F6").Value = 85
$ws.Range("F7").Value = 9801
$ws.Range("F10").Value = 3910
$ws.Range("F14").Value = 210
$ws.Range("F17").Value = 55
$ws.Range("F20").Value = 630
$ws.Range("F21").Value = 3927
$ws.Range("F22").Value = 140
$ws.Range("F24").Value = 5399
$ws.Range("F26").Value = 2132
$ws.Range("F28").Value = 363
$ws.Range("F29").Value = 8052
$ws.Range("F31").Value = 8
$ws.Range("F32").Value = 2211
$ws.Range("F33").Value = 2221
$ws.Range("F34").Value = 1340
$ws.Range("F35").Value = 1320
$ws.Range("F37").Value = 30
$ws.Range("F38").Value = 279
$ws.Range("F44").Value = 1350
$ws.Range("F45").Value = 2129
$ws.Range("F46").Value = 138
$ws.Range("F47").Value = 232
$ws.Range("F48").Value = 1221
$ws.Range("F49").Value = 5

$ws = $wb.Worksheets.Item("演出")
$ws.Range("F4").Value = 150
$ws.Range("F5").Value = 2
$ws.Range("F11").Value = 127
$ws.Range("F20").Value = 19

$ws = $wb.Worksheets.Item("本地生活")
$ws.Range("F2").Value = 590
$ws.Range("F3").Value = 771
$ws.Range("F4").Value = 70

$ws = $wb.Worksheets.Item("全部类型")
$ws.Range("F3").Value = 771
$ws.Range("F4").Value = 70
$ws.Range("F6").Value = 5769
$ws.Range("F7").Value = 85
$ws.Range("F8").Value = 3910
$ws.Range("F14").Value = 55
$ws.Range("F16").Value = 150
$ws.Range("F18").Value = 630
$ws.Range("F19").Value = 3927
$ws.Range("F21").Value = 140
$ws.Range("F23").Value = 5399
$ws.Range("F25").Value = 2132
$ws.Range("F27").Value = 363
$ws.Range("F28").Value = 8052
$ws.Range("F30").Value = 2211
$ws.Range("F31").Value = 2221
$ws.Range("F32").Value = 1340
$ws.Range("F33").Value = 1320
$ws.Range("F34").Value = 30
$ws.Range("F35").Value = 279
$ws.Range("F42").Value = 1350
$ws.Range("F44").Value = 2129
$ws.Range("F45").Value = 138
$ws.Range("F46").Value = 232
$ws.Range("F48").Value = 19
$ws.Range("F49").Value = 1221
